# Updated symbol list on Sat Jan 21 21:27:11 UTC 2023 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# crypto-exchange-token rows on Sheet1. Values are written with a leading
# apostrophe so Excel stores them as literal text (matching the workbook's
# existing text-formatted Price/Volume columns) instead of re-interpreting
# them as numbers/percentages and rounding or reformatting the text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value  = "'303.65"
$ws.Range("E2").Value  = "'1.38%"

$ws.Range("D3").Value  = "'35.68"
$ws.Range("E3").Value  = "'10.81%"

$ws.Range("D4").Value  = "'5.056"

$ws.Range("D5").Value  = "'0.07805"
$ws.Range("E5").Value  = "'1.25%"

$ws.Range("D6").Value  = "'2.250"
$ws.Range("E6").Value  = "'-0.14%"

$ws.Range("D7").Value  = "'8.085"
$ws.Range("E7").Value  = "'1.97%"

$ws.Range("D8").Value  = "'4.047"
$ws.Range("E8").Value  = "'6.05%"

$ws.Range("D9").Value  = "'0.9303"
$ws.Range("E9").Value  = "'0.87%"

$ws.Range("D10").Value = "'0.09489"
$ws.Range("E10").Value = "'-4.34%"

$ws.Range("D11").Value = "'0.1822"
$ws.Range("E11").Value = "'3.18%"

$ws.Range("D12").Value = "'0.08549"
$ws.Range("E12").Value = "'1.46%"

$ws.Range("E13").Value = "'3.74%"

$ws.Range("D14").Value = "'0.09945"
$ws.Range("E14").Value = "'1.12%"

$ws.Range("D15").Value = "'0.001480"
$ws.Range("E15").Value = "'0.06%"

$ws.Range("D16").Value = "'0.005748"
$ws.Range("E16").Value = "'1.67%"

$ws.Range("E17").Value = "'-1.54%"

$ws.Range("E18").Value = "'-0.79%"

$ws.Range("E19").Value = "'0.99%"

$ws.Range("D20").Value = "'0.1321"

$ws.Range("D21").Value = "'4.564"
$ws.Range("E21").Value = "'10.79%"

$ws.Range("D22").Value = "'0.2237"
$ws.Range("E22").Value = "'7.14%"

$ws.Range("D23").Value = "'0.04680"
$ws.Range("E23").Value = "'3.36%"

$ws.Range("D24").Value = "'0.001243"
$ws.Range("E24").Value = "'2.31%"

$ws.Range("D25").Value = "'0.004541"
$ws.Range("E25").Value = "'3.95%"

$ws.Range("E26").Value = "'1.04%"

$ws.Range("E27").Value = "'-19.82%"

$ws.Range("D39").Value = "'0.01775"
$ws.Range("E39").Value = "'3.97%"

$ws.Range("D40").Value = "'0.04705"
$ws.Range("E40").Value = "'1.06%"

$ws.Range("D41").Value = "'0.007914"
$ws.Range("E41").Value = "'3.62%"

$ws.Range("D42").Value = "'0.1418"
$ws.Range("E42").Value = "'1.78%"

$ws.Range("D43").Value = "'0.008001"
$ws.Range("E43").Value = "'-18.00%"

$ws.Range("D44").Value = "'0.002220"
$ws.Range("E44").Value = "'6.49%"

$ws.Range("D45").Value = "'0.009102"
$ws.Range("E45").Value = "'-6.24%"

$ws.Range("D46").Value = "'0.00006192"
$ws.Range("E46").Value = "'2.20%"

$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.79%"

$ws.Range("D48").Value = "'4.055"
$ws.Range("E48").Value = "'45.12%"

$ws.Range("D49").Value = "'0.002692"
$ws.Range("E49").Value = "'35.66%"

$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.79%"

$ws.Range("D51").Value = "'0.0002001"
$ws.Range("E51").Value = "'0.79%"
